# Update the "Datos actualizados" timestamp cell (row 1, column A)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 9 de Agosto de 2020 a las 12:11"

# Refresh per-country COVID stats (columns: B=Casos totales, C=Nuevos casos,
# D=Casos activos, E=Recuperados, F=Casos criticos, G=Muertes hoy, H=Muertes)
# for the countries whose figures changed, without touching row order here.

# Banglades (row 18)
$ws.Cells.Item(18, 2).Value = 257600
$ws.Cells.Item(18, 3).Value = 2487
$ws.Cells.Item(18, 4).Value = 148370
$ws.Cells.Item(18, 5).Value = 105831
$ws.Cells.Item(18, 7).Value = 34
$ws.Cells.Item(18, 8).Value = 3399

# Indonesia (row 26)
$ws.Cells.Item(26, 2).Value = 125396
$ws.Cells.Item(26, 3).Value = 1893
$ws.Cells.Item(26, 4).Value = 80952
$ws.Cells.Item(26, 5).Value = 38721
$ws.Cells.Item(26, 7).Value = 65
$ws.Cells.Item(26, 8).Value = 5723

# Rumania (row 44)
$ws.Cells.Item(44, 2).Value = 61768
$ws.Cells.Item(44, 3).Value = 1145
$ws.Cells.Item(44, 4).Value = 30119
$ws.Cells.Item(44, 5).Value = 28949
$ws.Cells.Item(44, 7).Value = 41
$ws.Cells.Item(44, 8).Value = 2700

# Azerbaiyan (row 60)
$ws.Cells.Item(60, 2).Value = 33568
$ws.Cells.Item(60, 3).Value = 87
$ws.Cells.Item(60, 4).Value = 30364
$ws.Cells.Item(60, 5).Value = 2714
$ws.Cells.Item(60, 7).Value = 2
$ws.Cells.Item(60, 8).Value = 490

# El Salvador (row 73)
$ws.Cells.Item(73, 2).Value = 20423
$ws.Cells.Item(73, 3).Value = 445
$ws.Cells.Item(73, 4).Value = 9626
$ws.Cells.Item(73, 5).Value = 10248
$ws.Cells.Item(73, 7).Value = 13
$ws.Cells.Item(73, 8).Value = 549

# Consejo Danes para los Refugiados (row 87)
$ws.Cells.Item(87, 2).Value = 9454
$ws.Cells.Item(87, 3).Value = 18
$ws.Cells.Item(87, 4).Value = 8324
$ws.Cells.Item(87, 5).Value = 906
$ws.Cells.Item(87, 7).Value = 6
$ws.Cells.Item(87, 8).Value = 224

# Malasia (row 88)
$ws.Cells.Item(88, 2).Value = 9083
$ws.Cells.Item(88, 3).Value = 13
$ws.Cells.Item(88, 4).Value = 8784
$ws.Cells.Item(88, 5).Value = 174

# Finlandia (row 95)
$ws.Cells.Item(95, 2).Value = 7584
$ws.Cells.Item(95, 3).Value = 16
$ws.Cells.Item(95, 5).Value = 273

# Eslovaquia (row 123)
$ws.Cells.Item(123, 2).Value = 2596
$ws.Cells.Item(123, 3).Value = 30
$ws.Cells.Item(123, 4).Value = 1864
$ws.Cells.Item(123, 5).Value = 701

# The small group Surinam/Eslovenia/Mozambique/Lituania/Estonia (rows 126-130)
# gets re-sorted because Lituania's new total (2252) now outranks Eslovenia
# (2247) and Mozambique (2241); Surinam (row 126, 2306) keeps its place.
# Row 127 becomes Lituania, row 128 becomes Eslovenia, row 129 becomes
# Mozambique (their own totals are unchanged), row 130 stays Estonia with an
# updated total.

# Row 127 -> Lituania (new figures)
$ws.Cells.Item(127, 1).Value = "Lituania"
$ws.Cells.Item(127, 2).Value = 2252
$ws.Cells.Item(127, 3).Value = 21
$ws.Cells.Item(127, 4).Value = 1670
$ws.Cells.Item(127, 5).Value = 501
$ws.Cells.Item(127, 8).Value = 81

# Row 128 -> Eslovenia (figures unchanged, just shifted down one row)
$ws.Cells.Item(128, 1).Value = "Eslovenia"
$ws.Cells.Item(128, 2).Value = 2247
$ws.Cells.Item(128, 3).Value = 0
$ws.Cells.Item(128, 4).Value = 1927
$ws.Cells.Item(128, 5).Value = 194
$ws.Cells.Item(128, 8).Value = 126

# Row 129 -> Mozambique (figures unchanged, just shifted down one row)
$ws.Cells.Item(129, 1).Value = "Mozambique"
$ws.Cells.Item(129, 2).Value = 2241
$ws.Cells.Item(129, 3).Value = 0
$ws.Cells.Item(129, 4).Value = 832
$ws.Cells.Item(129, 5).Value = 1393
$ws.Cells.Item(129, 8).Value = 16

# Row 130 -> Estonia (new figures, stays in place)
$ws.Cells.Item(130, 2).Value = 2152
$ws.Cells.Item(130, 3).Value = 5
$ws.Cells.Item(130, 5).Value = 128
